# Shane Watson vs Chennai Super Kings — refresh innings-by-innings batting log
# "updated activity till excel form": the existing 10 innings rows are
# reshuffled into a new order and a new 11th innings (row 6: 8 runs off 3
# balls, 2 fours, 0 sixes) is inserted, extending the table to row 12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final runs/balls/fours/sixes values for data rows 2..12 (in sheet order).
$data = @(
  @("14","19","1","1"),
  @("83","53","11","3"),
  @("50","40","6","1"),
  @("33","21","1","4"),
  @("8","3","2","0"),
  @("42","38","1","3"),
  @("36","28","6","0"),
  @("1","6","0","0"),
  @("4","5","1","0"),
  @("14","18","3","0"),
  @("14","16","1","1")
)

$firstRow = 2
$lastRow = $firstRow + $data.Length - 1

# New row 12 needs playerName / teamName filled in too (it didn't exist before).
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "Shane Watson "
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "Chennai Super Kings"

# Keep these columns stored as text (matches t="str" / numberStoredAsText
# semantics already used throughout the sheet) while writing the new values.
$ws.Range("C" + $firstRow + ":F" + $lastRow).NumberFormat = "@"

for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $firstRow + $i
  $vals = $data[$i]
  $ws.Cells.Item($row, 3).Value = $vals[0]
  $ws.Cells.Item($row, 4).Value = $vals[1]
  $ws.Cells.Item($row, 5).Value = $vals[2]
  $ws.Cells.Item($row, 6).Value = $vals[3]
}
